$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns updated with refreshed crypto data.
# Some D values (e.g. "1.003") are numeric-looking, so a leading apostrophe
# forces Excel to store them as literal text, matching the original inline
# string cells (e.g. "28.751.42" stays text because of the double dot,
# while "1.003" would otherwise be parsed as the number 1.003).
$ws.Range('D2').Value = '28.718.73'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '1.803.15'
$ws.Range('E3').Value = '  -1.26%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  +0.40%  '
$ws.Range('D5').Value = '''231.66'
$ws.Range('E5').Value = '  -2.03%  '
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').Value = '''0.2779'
$ws.Range('E8').Value = '  -1.47%  '
$ws.Range('D9').Value = '''0.06834'
$ws.Range('E9').Value = '  -3.90%  '
$ws.Range('D10').Value = '''23.34'
$ws.Range('E10').Value = '  -1.38%  '
$ws.Range('D11').Value = '''0.07541'
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').Value = '1.807.88'
$ws.Range('E12').Value = '  -0.83%  '
$ws.Range('D13').Value = '''4.800'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '''0.6252'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('D15').Value = '2.048.66'
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').Value = '''0.000009272'
$ws.Range('E16').Value = '  -8.26%  '
$ws.Range('E17').Value = '  -4.70%  '
$ws.Range('D18').Value = '28.700.44'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('E19').Value = '  -6.93%  '
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').Value = '''210.66'
$ws.Range('E21').Value = '  -7.61%  '
$ws.Range('D22').Value = '''11.46'
$ws.Range('E22').Value = '  -2.87%  '
$ws.Range('D23').Value = '''6.857'
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('D24').Value = '''1.003'
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('D25').Value = '''154.26'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('D26').Value = '''7.842'
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('D27').Value = '''0.1277'
$ws.Range('E27').Value = '  -3.57%  '
$ws.Range('D28').Value = '''16.41'
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('D29').Value = '''1.433'
$ws.Range('E29').Value = '  -3.63%  '
$ws.Range('D30').Value = '''0.06196'
$ws.Range('E30').Value = '  -2.71%  '
$ws.Range('E31').Value = '  -2.18%  '
$ws.Range('D32').Value = '''3.780'
$ws.Range('E32').Value = '  -1.25%  '
$ws.Range('D33').Value = '''3.749'
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('D34').Value = '''1.718'
$ws.Range('E34').Value = '  -1.82%  '
$ws.Range('D35').Value = '''1.061'
$ws.Range('E35').Value = '  -5.81%  '
$ws.Range('D36').Value = '''0.6403'
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('E37').Value = '  -1.93%  '
$ws.Range('D38').Value = '''2.719'
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').Value = '''0.01712'
$ws.Range('E39').Value = '  -1.78%  '
$ws.Range('D40').Value = '''6.431'
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('D41').Value = '1.134.42'
$ws.Range('E41').Value = '  -6.92%  '
$ws.Range('D42').Value = '''0.8678'
$ws.Range('E42').Value = '  -6.52%  '
$ws.Range('D43').Value = '''1.003'
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('D44').Value = '''100.58'
$ws.Range('E44').Value = '  -0.70%  '
$ws.Range('D45').Value = '1.965.97'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('D46').Value = '''60.59'
$ws.Range('E46').Value = '  -3.96%  '
$ws.Range('E47').Value = '  -5.19%  '
$ws.Range('D48').Value = '''1.597'
$ws.Range('E48').Value = '  -1.89%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.05474'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''8.333'
$ws.Range('E50').Value = '  -3.37%  '
$ws.Range('E51').Value = '  -1.49%  '
